$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 ("I0") and J1 ("IF") with the same formatting as
# the existing header cells (bold font + border + centered alignment).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$iValues = @{
    2  = 1
    3  = 1
    4  = 1
    5  = 1
    6  = 1
    7  = 1
    8  = 1
    9  = 1
    10 = 1
    11 = 1
    12 = 1
    13 = 1
    14 = 1
    15 = 1
    16 = 1
    17 = 1
    18 = 1
    19 = 7
    20 = 9
    21 = 4
    22 = 4
    23 = 8
}

$jValues = @{
    2  = 6
    3  = 5
    4  = 6
    5  = 6
    6  = 6
    7  = 6
    8  = 5
    9  = 2
    10 = 7
    11 = 6
    12 = 6
    13 = 4
    14 = 4
    15 = 6
    16 = 6
    17 = 6
    18 = 7
    19 = 9
    20 = 9
    21 = 6
    22 = 5
    23 = 9
}

for ($row = 2; $row -le 23; $row++) {
    $ws.Cells.Item($row, 9).Value = $iValues[$row]
    $ws.Cells.Item($row, 10).Value = $jValues[$row]
}
